$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 14).ClearContents()
$ws.Cells.Item(40, 8).Value = 2300.1
$ws.Cells.Item(40, 9).Value = 1500
$ws.Cells.Item(40, 11).Value = 1500
$ws.Cells.Item(40, 13).Value = -1325
$ws.Cells.Item(58, 8).Value = 338
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 14).ClearContents()
$ws.Cells.Item(80, 8).Value = 1652.6
$ws.Cells.Item(80, 10).Value = 2207.8
$ws.Cells.Item(80, 12).Value = 6623.400000000001
$ws.Cells.Item(80, 14).Value = -8619.400000000001
$ws.Cells.Item(83, 8).Value = 1652.6
$ws.Cells.Item(83, 10).Value = 2207.8
$ws.Cells.Item(83, 12).Value = 19870.2
$ws.Cells.Item(83, 14).Value = -29854.2
$ws.Cells.Item(113, 8).Value = 8114.4614
$ws.Cells.Item(113, 10).Value = 8221
$ws.Cells.Item(113, 12).Value = 8221
$ws.Cells.Item(113, 14).Value = -14729
$ws.Cells.Item(138, 8).Value = 2634.7273
$ws.Cells.Item(138, 9).Value = 1995.5
$ws.Cells.Item(138, 11).Value = 5986.5
$ws.Cells.Item(138, 13).Value = -846.5
$ws.Cells.Item(141, 8).Value = 2283.2856
$ws.Cells.Item(141, 9).Value = 2283.2856
$ws.Cells.Item(141, 11).Value = 6849.8568
$ws.Cells.Item(141, 13).Value = -1669.8568
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4010.7693
$ws.Cells.Item(32, 9).Value = 1567.2778
$ws.Cells.Item(32, 11).Value = 1567.2778
$ws.Cells.Item(32, 13).Value = -1280.2778
$ws.Cells.Item(62, 8).Value = 56742
$ws.Cells.Item(62, 9).Value = 60226
$ws.Cells.Item(62, 10).Value = 55000
$ws.Cells.Item(62, 11).Value = 60226
$ws.Cells.Item(62, 12).Value = 55000
$ws.Cells.Item(62, 13).Value = -59602
$ws.Cells.Item(62, 14).Value = -56248
$ws.Cells.Item(65, 8).Value = 56742
$ws.Cells.Item(65, 9).Value = 60226
$ws.Cells.Item(65, 10).Value = 55000
$ws.Cells.Item(65, 11).Value = 180678
$ws.Cells.Item(65, 12).Value = 165000
$ws.Cells.Item(65, 13).Value = -177558
$ws.Cells.Item(65, 14).Value = -171240
$ws.Cells.Item(102, 8).Value = 2900
$ws.Cells.Item(102, 9).Value = 2900
$ws.Cells.Item(102, 11).Value = 2900
$ws.Cells.Item(102, 13).Value = -1278
$ws.Cells.Item(122, 8).Value = 3149.1667
$ws.Cells.Item(122, 9).Value = 2979
$ws.Cells.Item(122, 11).Value = 8937
$ws.Cells.Item(122, 13).Value = -6487
$ws.Cells.Item(133, 8).Value = 80000
$ws.Cells.Item(133, 10).Value = 80000
$ws.Cells.Item(133, 12).Value = 80000
$ws.Cells.Item(133, 14).Value = -85060
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 14).ClearContents()
$ws.Cells.Item(19, 8).Value = 25798.092
$ws.Cells.Item(19, 10).Value = 24000
$ws.Cells.Item(19, 12).Value = 24000
$ws.Cells.Item(19, 14).Value = -24346
$ws.Cells.Item(20, 8).Value = 3699.818
$ws.Cells.Item(20, 9).Value = 4462.25
$ws.Cells.Item(20, 10).Value = 1666.6666
$ws.Cells.Item(20, 11).Value = 4462.25
$ws.Cells.Item(20, 12).Value = 1666.6666
$ws.Cells.Item(20, 13).Value = -4215.25
$ws.Cells.Item(20, 14).Value = -2160.6666
$ws.Cells.Item(134, 8).Value = 12400
$ws.Cells.Item(134, 9).Value = 12400
$ws.Cells.Item(134, 11).Value = 37200
$ws.Cells.Item(134, 13).Value = -34665
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2251.0417
$ws.Cells.Item(31, 9).Value = 1815.9166
$ws.Cells.Item(31, 11).Value = 1815.9166
$ws.Cells.Item(31, 13).Value = -1520.9166
$ws.Cells.Item(34, 8).Value = 2251.0417
$ws.Cells.Item(34, 9).Value = 1815.9166
$ws.Cells.Item(34, 11).Value = 1815.9166
$ws.Cells.Item(34, 13).Value = -1613.9166
$ws.Cells.Item(41, 8).Value = 10722.429
$ws.Cells.Item(41, 10).Value = 12499.667
$ws.Cells.Item(41, 12).Value = 12499.667
$ws.Cells.Item(41, 14).Value = -13355.667
$ws.Cells.Item(50, 8).Value = 20216.6
$ws.Cells.Item(59, 8).Value = 29084
$ws.Cells.Item(60, 8).Value = 22015.166
$ws.Cells.Item(105, 8).Value = 595
$ws.Cells.Item(105, 9).Value = 595
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 595
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = 1152
$ws.Cells.Item(105, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 5170.857
$ws.Cells.Item(132, 9).Value = 5049
$ws.Cells.Item(132, 11).Value = 15147
$ws.Cells.Item(132, 13).Value = -12617
$ws.Cells.Item(134, 8).Value = 5000
$ws.Cells.Item(134, 9).Value = 5000
$ws.Cells.Item(134, 10).Value = 5000
$ws.Cells.Item(134, 11).Value = 15000
$ws.Cells.Item(134, 12).Value = 15000
$ws.Cells.Item(134, 13).Value = -12465
$ws.Cells.Item(134, 14).Value = -20070
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1749
$ws.Cells.Item(5, 9).Value = 1248.5
$ws.Cells.Item(5, 11).Value = 3745.5
$ws.Cells.Item(5, 13).Value = -3633.5
$ws.Cells.Item(6, 8).Value = 2303.4666
$ws.Cells.Item(6, 9).Value = 333.84616
$ws.Cells.Item(6, 10).Value = 15106
$ws.Cells.Item(6, 11).Value = 1001.53848
$ws.Cells.Item(6, 12).Value = 45318
$ws.Cells.Item(6, 13).Value = -888.5384799999999
$ws.Cells.Item(6, 14).Value = -45544
$ws.Cells.Item(68, 8).Value = 948.5
$ws.Cells.Item(68, 9).Value = 899.5
$ws.Cells.Item(68, 10).Value = 997.5
$ws.Cells.Item(68, 11).Value = 2698.5
$ws.Cells.Item(68, 12).Value = 2992.5
$ws.Cells.Item(68, 13).Value = -1887.5
$ws.Cells.Item(68, 14).Value = -4614.5
$ws.Cells.Item(71, 8).Value = 948.5
$ws.Cells.Item(71, 9).Value = 899.5
$ws.Cells.Item(71, 10).Value = 997.5
$ws.Cells.Item(71, 11).Value = 8095.5
$ws.Cells.Item(71, 12).Value = 8977.5
$ws.Cells.Item(71, 13).Value = -4039.5
$ws.Cells.Item(71, 14).Value = -17089.5
$ws.Cells.Item(114, 8).Value = 650
$ws.Cells.Item(114, 9).Value = 650
$ws.Cells.Item(114, 10).Value = 0
$ws.Cells.Item(114, 11).Value = 1950
$ws.Cells.Item(114, 12).Value = 0
$ws.Cells.Item(114, 13).Value = 1304
$ws.Cells.Item(114, 14).ClearContents()
$ws.Cells.Item(129, 8).Value = 1734
$ws.Cells.Item(129, 9).Value = 1005
$ws.Cells.Item(129, 11).Value = 3015
$ws.Cells.Item(129, 13).Value = 1985
$ws.Cells.Item(135, 8).Value = 1749
$ws.Cells.Item(135, 9).Value = 1248.5
$ws.Cells.Item(135, 11).Value = 11236.5
$ws.Cells.Item(135, 13).Value = -8701.5
$ws.Cells.Item(137, 8).Value = 2029
$ws.Cells.Item(137, 9).Value = 2029
$ws.Cells.Item(137, 11).Value = 6087
$ws.Cells.Item(137, 13).Value = -987
$ws.Cells.Item(140, 8).Value = 2484.125
$ws.Cells.Item(140, 9).Value = 1982.2858
$ws.Cells.Item(140, 11).Value = 5946.857400000001
$ws.Cells.Item(140, 13).Value = -766.8574000000008
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(46, 8).Value = 18304.1
$ws.Cells.Item(46, 9).Value = 13260.25
$ws.Cells.Item(46, 10).Value = 21666.666
$ws.Cells.Item(46, 11).Value = 13260.25
$ws.Cells.Item(46, 12).Value = 21666.666
$ws.Cells.Item(46, 13).Value = -13104.25
$ws.Cells.Item(46, 14).Value = -21978.666
$ws.Cells.Item(97, 8).Value = 3161.3
$ws.Cells.Item(97, 9).Value = 2326.75
$ws.Cells.Item(97, 10).Value = 6499.5
$ws.Cells.Item(97, 11).Value = 2326.75
$ws.Cells.Item(97, 12).Value = 6499.5
$ws.Cells.Item(97, 13).Value = -1830.75
$ws.Cells.Item(97, 14).Value = -7491.5
$ws.Cells.Item(123, 8).Value = 50000
$ws.Cells.Item(123, 10).Value = 50000
$ws.Cells.Item(123, 12).Value = 50000
$ws.Cells.Item(123, 14).Value = -54900
$ws.Cells.Item(132, 8).Value = 6612
$ws.Cells.Item(132, 9).Value = 6612
$ws.Cells.Item(132, 11).Value = 19836
$ws.Cells.Item(132, 13).Value = -17306
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(123, 8).Value = 15000
$ws.Cells.Item(123, 9).Value = 15000
$ws.Cells.Item(123, 11).Value = 15000
$ws.Cells.Item(123, 13).Value = -10100
$ws.Cells.Item(132, 8).Value = 2500
$ws.Cells.Item(132, 9).Value = 2500
$ws.Cells.Item(132, 11).Value = 7500
$ws.Cells.Item(132, 13).Value = -4970
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 1306.9166
$ws.Cells.Item(107, 9).Value = 910.375
$ws.Cells.Item(107, 11).Value = 2731.125
$ws.Cells.Item(107, 13).Value = -811.125
$ws.Cells.Item(132, 8).Value = 3116.875
$ws.Cells.Item(132, 9).Value = 1988
$ws.Cells.Item(132, 10).Value = 4998.3335
$ws.Cells.Item(132, 11).Value = 5964
$ws.Cells.Item(132, 12).Value = 14995.0005
$ws.Cells.Item(132, 13).Value = -3434
$ws.Cells.Item(132, 14).Value = -20055.0005
